$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and Report Covering Week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/17/2023  Through  7/23/2023"

# --- Crime statistics table updates (rows 14-29) ---
# Row 14
$ws.Range("F14").Value = 1

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 4
$ws.Range("I15").Value = 18
$ws.Range("K15").Value = 12.5
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 63.636363636363
$ws.Range("N15").Value = 5.882352941176

# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = -35.714285714285
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 35
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 174
$ws.Range("J16").Value = 164
$ws.Range("K16").Value = 6.097560975609
$ws.Range("L16").Value = 75.757575757575
$ws.Range("M16").Value = 21.678321678321
$ws.Range("N16").Value = -62.58064516129

# Row 17
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 183.333333333333
$ws.Range("F17").Value = 60
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = 57.894736842105
$ws.Range("I17").Value = 262
$ws.Range("J17").Value = 232
$ws.Range("K17").Value = 12.931034482758
$ws.Range("L17").Value = 48.863636363636
$ws.Range("M17").Value = 20.183486238532
$ws.Range("N17").Value = -0.380228136882

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 102
$ws.Range("J18").Value = 85
$ws.Range("K18").Value = 20
$ws.Range("L18").Value = 108.163265306122
$ws.Range("M18").Value = -14.285714285714
$ws.Range("N18").Value = -81.386861313868

# Row 19
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 71.428571428571
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = -9.090909090909
$ws.Range("I19").Value = 230
$ws.Range("J19").Value = 244
$ws.Range("K19").Value = -5.737704918032
$ws.Range("L19").Value = 35.294117647058
$ws.Range("M19").Value = 54.362416107382
$ws.Range("N19").Value = -6.122448979591

# Row 20
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 60
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 63.157894736842
$ws.Range("I20").Value = 189
$ws.Range("J20").Value = 131
$ws.Range("K20").Value = 44.274809160305
$ws.Range("L20").Value = 139.240506329114
$ws.Range("M20").Value = 350
$ws.Range("N20").Value = -16.740088105726

# Row 21
$ws.Range("C21").Value = 53
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 47.222222222222
$ws.Range("F21").Value = 176
$ws.Range("G21").Value = 145
$ws.Range("H21").Value = 21.379310344827
$ws.Range("I21").Value = 981
$ws.Range("J21").Value = 875
$ws.Range("K21").Value = 12.114285714285
$ws.Range("L21").Value = 63.772954924874
$ws.Range("M21").Value = 43.421052631578
$ws.Range("N21").Value = -45.256696428571

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 13
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = -23.529411764705
$ws.Range("L22").Value = 8.333333333333
$ws.Range("M22").Value = 0

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 2
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = -50
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = 12
$ws.Range("K23").Value = 25
$ws.Range("L23").Value = 7.142857142857
$ws.Range("M23").Value = 114.285714285714

# Row 24
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 220
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 55
$ws.Range("I24").Value = 467
$ws.Range("J24").Value = 460
$ws.Range("K24").Value = 1.521739130434
$ws.Range("L24").Value = 72.324723247232
$ws.Range("M24").Value = 47.318611987381

# Row 25
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 12.5
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 61
$ws.Range("H25").Value = -13.11475409836
$ws.Range("I25").Value = 299
$ws.Range("J25").Value = 325
$ws.Range("K25").Value = -8
$ws.Range("L25").Value = 43.75
$ws.Range("M25").Value = -17.403314917127

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("F26").Value = 7
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 600
$ws.Range("I26").Value = 34
$ws.Range("K26").Value = 17.241379310344
$ws.Range("L26").Value = -17.073170731707

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 56
$ws.Range("J27").Value = 63
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = -16.417910447761

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 200
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = 7.142857142857
$ws.Range("M28").Value = -6.25
$ws.Range("N28").Value = -72.727272727272

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("J29").Value = 12
$ws.Range("K29").Value = 8.333333333333
$ws.Range("M29").Value = -7.142857142857
$ws.Range("N29").Value = -73.469387755102
